$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.775.00"
$ws.Range("E2").Value = "  +3.54%  "
$ws.Range("D3").Value = "3.625.13"
$ws.Range("E3").Value = "  +6.49%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "593.00"
$ws.Range("E5").Value = "  +0.93%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "180.42"
$ws.Range("E6").Value = "  -0.65%  "
$ws.Range("E7").Value = "  +6.46%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.612"
$ws.Range("E8").Value = "  +2.14%  "
$ws.Range("E9").Value = "  +0.19%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.203"
$ws.Range("E10").Value = "  +2.20%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.606"
$ws.Range("E11").Value = "  +1.95%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "49.93"
$ws.Range("E12").Value = "  +2.69%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000286"
$ws.Range("E13").Value = "  +0.57%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "696.75"
$ws.Range("E14").Value = "  +2.02%  "
$ws.Range("D15").Value = "4.217.63"
$ws.Range("E15").Value = "  +6.57%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.99"
$ws.Range("E16").Value = "  +3.38%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.666.98"
$ws.Range("E17").Value = "  +7.80%  "
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "71.918.59"
$ws.Range("E18").Value = "  +3.59%  "
$ws.Range("E19").Value = "  +1.99%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.39"
$ws.Range("E20").Value = "  +3.59%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.63"
$ws.Range("E21").Value = "  +2.51%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.936"
$ws.Range("E22").Value = "  +2.70%  "
$ws.Range("E23").Value = "  +8.38%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "17.89"
$ws.Range("E24").Value = "  +3.36%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "103.73"
$ws.Range("E25").Value = "  +0.27%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.03"
$ws.Range("E26").Value = "  +2.25%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.86"
$ws.Range("E27").Value = "  +4.61%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.03"
$ws.Range("E28").Value = "  +3.33%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "35.06"
$ws.Range("E29").Value = "  +3.09%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.18"
$ws.Range("E30").Value = "  +4.18%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.31"
$ws.Range("E31").Value = "  +4.71%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.18"
$ws.Range("E32").Value = "  +14.78%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "583.88"
$ws.Range("E33").Value = "  +4.15%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.36"
$ws.Range("E34").Value = "  +1.67%  "
$ws.Range("E35").Value = "  +2.70%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "59.55"
$ws.Range("E36").Value = "  +1.34%  "
$ws.Range("D38").Value = "3.662.93"
$ws.Range("E38").Value = "  -0.07%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.143"
$ws.Range("E39").Value = "  +1.20%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "35.92"
$ws.Range("E40").Value = "  -0.39%  "
$ws.Range("D41").Value = "0.0₃0765"
$ws.Range("E41").Value = "  +5.24%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.42"
$ws.Range("E42").Value = "  +4.60%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0466"
$ws.Range("E43").Value = "  +9.13%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.76"
$ws.Range("E44").Value = "  +2.88%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.350"
$ws.Range("E45").Value = "  +2.88%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.41"
$ws.Range("E46").Value = "  +2.51%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.82"
$ws.Range("E47").Value = "  +5.33%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.133"
$ws.Range("E48").Value = "  +2.31%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.44"
$ws.Range("E49").Value = "  +3.57%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.00"
$ws.Range("E50").Value = "  -0.18%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "131.78"
$ws.Range("E51").Value = "  -0.67%  "
